$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Kagiso Rabada"

# Force the whole used range to be stored as text (matches source data
# which keeps numeric-looking values like "9", "225.00" as strings).
$ws.Range("A1:M4").NumberFormat = "@"

# Header row
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - new match (41st)
$row2 = @("41st","Delhi Capitals","Kagiso Rabada","","0","1","0","0","0.00","Kolkata Knight Riders","Sharjah","September 28","KKR won by 3 wickets (with 10 balls remaining)")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# Row 3 - new match (50th)
$row3 = @("50th","Delhi Capitals","Kagiso Rabada","","4","1","1","0","400.00","Chennai Super Kings","Dubai (DSC)","October 04","Capitals won by 3 wickets (with 2 balls remaining)")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Row 4 - original match data (7th), shifted down from row 2 and right by one column
$row4 = @("7th","Delhi Capitals","Kagiso Rabada","","9","4","1","0","225.00","Rajasthan Royals","Wankhede","April 15","Royals won by 3 wickets (with 2 balls remaining)")
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}
